$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(80, 1).Value = 78
$ws.Cells.Item(80, 2).Value = "Okulun ücretli-ücretsiz shuttlelarında neden hiç geç saatlere kampüse ulaşım sağlayamıyoruz, sağlayabilecek bir sefer yok? Öğrenciyiz, yurtlarda kalıyoruz Şile'ye genel olarak ulaşım sıkıntılı geç saatlerde (00:00 vs.) ulaşılabilecek hiçbir olasılığımız yok. Biz bir konsere bir tiyatroya akşam vakti gidemeyecek miyiz, gittik nasıl döneceğiz?"
$ws.Cells.Item(81, 1).Value = 79
$ws.Cells.Item(81, 2).Value = "İnternet alt yapısının sıkıntıları yüzünden derslerime giremiyorum. Dünden beri WeTransfer'den dosya göndermeye çalışıyorum. 1.3 gb'lık dosyayı 24 saat oldu şuan 12 saat daha ekledi."
$ws.Cells.Item(82, 1).Value = 80
$ws.Cells.Item(82, 2).Value = "Selamlar, kampüse gelmeden online olarak ilişiği kesebiliyor muyuz?"
$ws.Cells.Item(83, 1).Value = 81
$ws.Cells.Item(83, 2).Value = "Okulun satranç klübü var mı?"
$ws.Cells.Item(84, 1).Value = 82
$ws.Cells.Item(84, 2).Value = "İnternete bağlanmakta sorun yaşıyorum iki gündür. Beyaz ekranda kalıyor. Online derslerime bağlanamıyorum. Bu internet sorunları ne zaman düzelecek Işık Üniversitesi?"
$ws.Cells.Item(85, 1).Value = 83
$ws.Cells.Item(85, 2).Value = "DMF'deki dolapları direkt olarak kullanabiliyor muyuz yoksa öncesinde dilekçe, ücret vs. gerekiyor mu?"
$ws.Cells.Item(86, 1).Value = 84
$ws.Cells.Item(86, 2).Value = "Müzik evini kullanmak serbest mi acaba? Randevu falan alınıyor mu?"
$ws.Cells.Item(87, 1).Value = 85
$ws.Cells.Item(87, 2).Value = "Kampüs içinde olan Komogene'nin numarasını bilen varsa yazabilir mi?"
$ws.Cells.Item(88, 1).Value = 86
$ws.Cells.Item(88, 2).Value = "Shuttle Kozyatağı'nda nereden geçiyor?"
$ws.Cells.Item(89, 1).Value = 87
$ws.Cells.Item(89, 2).Value = "Merhabalar, yurtların alt katındaki misafirhaneler hakkında bilgisi olan var mı?"
$ws.Cells.Item(90, 1).Value = 88
$ws.Cells.Item(90, 2).Value = "Merhaba, Türkçe bölüm okuyorum, İngilizce bir ders seçmiştim. O dersin Türkçe'sine nasıl geçiş yapabilirim, bilginiz var mı?"
$ws.Cells.Item(91, 1).Value = 89
$ws.Cells.Item(91, 2).Value = "Kadıköy'den kalkan shuttle Atatürk Caddesi Starbucks'ın önünde duruyor mu?"
$ws.Cells.Item(92, 1).Value = 90
$ws.Cells.Item(92, 2).Value = "Kadıköy'den kalkan shuttle tam olarak nerede ya Libadiye'den nasıl gidebilirim?"

$ws.Cells.Item(95, 4).Select()

